$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = "D2"; Value = "328.05"; ForceText = 1 }
    @{ Ref = "E2"; Value = "-0.19%"; ForceText = 1 }
    @{ Ref = "D3"; Value = "44.27"; ForceText = 1 }
    @{ Ref = "E3"; Value = "-0.25%"; ForceText = 1 }
    @{ Ref = "D4"; Value = "5.566"; ForceText = 1 }
    @{ Ref = "E4"; Value = "1.56%"; ForceText = 1 }
    @{ Ref = "D5"; Value = "0.08060"; ForceText = 1 }
    @{ Ref = "E5"; Value = "-1.33%"; ForceText = 1 }
    @{ Ref = "D6"; Value = "1.910"; ForceText = 1 }
    @{ Ref = "E6"; Value = "0.14%"; ForceText = 1 }
    @{ Ref = "D8"; Value = "0.9506"; ForceText = 1 }
    @{ Ref = "E8"; Value = "0.98%"; ForceText = 1 }
    @{ Ref = "D9"; Value = "0.1213"; ForceText = 1 }
    @{ Ref = "E9"; Value = "1.50%"; ForceText = 1 }
    @{ Ref = "D10"; Value = "0.1847"; ForceText = 1 }
    @{ Ref = "E10"; Value = "-2.86%"; ForceText = 1 }
    @{ Ref = "B11"; Value = "MCDex"; ForceText = 0 }
    @{ Ref = "C11"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"; ForceText = 0 }
    @{ Ref = "D11"; Value = "9.998"; ForceText = 1 }
    @{ Ref = "E11"; Value = "13.74%"; ForceText = 1 }
    @{ Ref = "B12"; Value = "MandalaExchangeToken"; ForceText = 0 }
    @{ Ref = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; ForceText = 0 }
    @{ Ref = "D12"; Value = "0.09649"; ForceText = 1 }
    @{ Ref = "E12"; Value = "-2.28%"; ForceText = 1 }
    @{ Ref = "B13"; Value = "BitrueCoin"; ForceText = 0 }
    @{ Ref = "C13"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; ForceText = 0 }
    @{ Ref = "D13"; Value = "0.04364"; ForceText = 1 }
    @{ Ref = "E13"; Value = "4.23%"; ForceText = 1 }
    @{ Ref = "B14"; Value = "BitMartToken"; ForceText = 0 }
    @{ Ref = "C14"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; ForceText = 0 }
    @{ Ref = "D14"; Value = "0.1064"; ForceText = 1 }
    @{ Ref = "E14"; Value = "-0.42%"; ForceText = 1 }
    @{ Ref = "B15"; Value = "BitForexToken"; ForceText = 0 }
    @{ Ref = "C15"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; ForceText = 0 }
    @{ Ref = "D15"; Value = "0.001286"; ForceText = 1 }
    @{ Ref = "E15"; Value = "-0.34%"; ForceText = 1 }
    @{ Ref = "B16"; Value = "CoinExToken"; ForceText = 0 }
    @{ Ref = "C16"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; ForceText = 0 }
    @{ Ref = "D16"; Value = "0.04214"; ForceText = 1 }
    @{ Ref = "E16"; Value = "-3.77%"; ForceText = 1 }
    @{ Ref = "B17"; Value = "TigerCash"; ForceText = 0 }
    @{ Ref = "C17"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; ForceText = 0 }
    @{ Ref = "D17"; Value = "0.005944"; ForceText = 1 }
    @{ Ref = "E17"; Value = "-2.31%"; ForceText = 1 }
    @{ Ref = "B18"; Value = "LEO"; ForceText = 0 }
    @{ Ref = "C18"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; ForceText = 0 }
    @{ Ref = "D18"; Value = "3.394"; ForceText = 1 }
    @{ Ref = "E18"; Value = "-4.03%"; ForceText = 1 }
    @{ Ref = "B19"; Value = "GateToken"; ForceText = 0 }
    @{ Ref = "C19"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; ForceText = 0 }
    @{ Ref = "D19"; Value = "4.282"; ForceText = 1 }
    @{ Ref = "E19"; Value = "-0.98%"; ForceText = 1 }
    @{ Ref = "B20"; Value = "BitpandaEcosystemToken"; ForceText = 0 }
    @{ Ref = "C20"; Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; ForceText = 0 }
    @{ Ref = "D20"; Value = "0.3465"; ForceText = 1 }
    @{ Ref = "E20"; Value = "-1.05%"; ForceText = 1 }
    @{ Ref = "D21"; Value = "0.1423"; ForceText = 1 }
    @{ Ref = "E21"; Value = "5.36%"; ForceText = 1 }
    @{ Ref = "D22"; Value = "0.2503"; ForceText = 1 }
    @{ Ref = "E22"; Value = "0.22%"; ForceText = 1 }
    @{ Ref = "D23"; Value = "0.001244"; ForceText = 1 }
    @{ Ref = "E23"; Value = "0.34%"; ForceText = 1 }
    @{ Ref = "D24"; Value = "0.004304"; ForceText = 1 }
    @{ Ref = "E24"; Value = "0.03%"; ForceText = 1 }
    @{ Ref = "D25"; Value = "0.0001191"; ForceText = 1 }
    @{ Ref = "E25"; Value = "-3.61%"; ForceText = 1 }
    @{ Ref = "E26"; Value = "-0.73%"; ForceText = 1 }
    @{ Ref = "D38"; Value = "0.02679"; ForceText = 1 }
    @{ Ref = "D39"; Value = "0.05521"; ForceText = 1 }
    @{ Ref = "E39"; Value = "-3.18%"; ForceText = 1 }
    @{ Ref = "D40"; Value = "0.007606"; ForceText = 1 }
    @{ Ref = "E40"; Value = "-3.51%"; ForceText = 1 }
    @{ Ref = "D41"; Value = "0.1404"; ForceText = 1 }
    @{ Ref = "E41"; Value = "-0.58%"; ForceText = 1 }
    @{ Ref = "D42"; Value = "0.007931"; ForceText = 1 }
    @{ Ref = "E42"; Value = "-18.54%"; ForceText = 1 }
    @{ Ref = "D43"; Value = "0.002018"; ForceText = 1 }
    @{ Ref = "E43"; Value = "-4.16%"; ForceText = 1 }
    @{ Ref = "D44"; Value = "0.008895"; ForceText = 1 }
    @{ Ref = "E44"; Value = "-8.07%"; ForceText = 1 }
    @{ Ref = "D45"; Value = "0.00007111"; ForceText = 1 }
    @{ Ref = "E45"; Value = "0.81%"; ForceText = 1 }
    @{ Ref = "E46"; Value = "-0.40%"; ForceText = 1 }
    @{ Ref = "D47"; Value = "0.002842"; ForceText = 1 }
    @{ Ref = "E47"; Value = "-17.48%"; ForceText = 1 }
    @{ Ref = "D48"; Value = "0.002268"; ForceText = 1 }
    @{ Ref = "E48"; Value = "-0.56%"; ForceText = 1 }
    @{ Ref = "E49"; Value = "-0.40%"; ForceText = 1 }
    @{ Ref = "E50"; Value = "-0.40%"; ForceText = 1 }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Ref)
    if ($u.ForceText -eq 1) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
